$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999438518516
$ws.Range("A2").Value = 0.99709017844233783
$ws.Range("A3").Value = 0.9895527881575763
$ws.Range("A4").Value = 0.98984830389634781
$ws.Range("A5").Value = 0.98039632986798819
$ws.Range("A6").Value = 0.9575811705487014
$ws.Range("A7").Value = 0.95261143209800347
$ws.Range("A8").Value = 0.94642387762550151
$ws.Range("A9").Value = 0.93949388314879378
$ws.Range("A10").Value = 0.93348392126411617
$ws.Range("A11").Value = 0.93268248310159363
$ws.Range("A12").Value = 0.93141341970622071
$ws.Range("A13").Value = 0.92698104950410376
$ws.Range("A14").Value = 0.9257344390380231
$ws.Range("A15").Value = 0.92596498182823184
$ws.Range("A16").Value = 0.92684935662824697
$ws.Range("A17").Value = 0.92314155078905069
$ws.Range("A18").Value = 0.92203263839515071
$ws.Range("A19").Value = 0.99420970493439964
$ws.Range("A20").Value = 0.97086382694604667
$ws.Range("A21").Value = 0.96457532149167524
$ws.Range("A22").Value = 0.96331080314314654
$ws.Range("A23").Value = 0.98251832884524182
$ws.Range("A24").Value = 0.96949779795199276
$ws.Range("A25").Value = 0.96304085794251781
$ws.Range("A26").Value = 0.95553072176953058
$ws.Range("A27").Value = 0.95292193479129761
$ws.Range("A28").Value = 0.94189552148916322
$ws.Range("A29").Value = 0.93438831134622835
$ws.Range("A30").Value = 0.93166737409809453
$ws.Range("A31").Value = 0.93391530115165211
$ws.Range("A32").Value = 0.93568645443940712
$ws.Range("A33").Value = 0.93516644213797351
